$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2065.5
$ws.Range("I38").Value = 2036.3334
$ws.Range("J38").Value = 2094.6667
$ws.Range("K38").Value = 6109.0002
$ws.Range("L38").Value = 6284.000100000001
$ws.Range("M38").Value = -5737.0002
$ws.Range("N38").Value = -7028.000100000001
$ws.Range("H101").Value = 9400
$ws.Range("I101").Value = 200
$ws.Range("J101").Value = 14000
$ws.Range("K101").Value = 600
$ws.Range("L101").Value = 42000
$ws.Range("M101").Value = 1022
$ws.Range("N101").Value = -45244
$ws.Range("H133").Value = 18371.25
$ws.Range("J133").Value = 18371.25
$ws.Range("L133").Value = 18371.25
$ws.Range("N133").Value = -28491.25
$ws.Range("H137").Value = 50001372
$ws.Range("I137").Value = 58824690
$ws.Range("J137").Value = 2567.6667
$ws.Range("K137").Value = 176474070
$ws.Range("L137").Value = 7703.000100000001
$ws.Range("M137").Value = -176471520
$ws.Range("N137").Value = -12803.0001
$ws.Range("H138").Value = 6582418.5
$ws.Range("I138").Value = 2980720.2
$ws.Range("J138").Value = 7755064.5
$ws.Range("K138").Value = 8942160.600000001
$ws.Range("L138").Value = 23265193.5
$ws.Range("M138").Value = -8937020.600000001
$ws.Range("N138").Value = -23275473.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18854.59
$ws.Range("I32").Value = 2533.76
$ws.Range("J32").Value = 93040.17999999999
$ws.Range("K32").Value = 2533.76
$ws.Range("L32").Value = 93040.17999999999
$ws.Range("M32").Value = -2246.76
$ws.Range("N32").Value = -93614.17999999999
$ws.Range("H45").Value = 898.1905
$ws.Range("I45").Value = 857.3333
$ws.Range("J45").Value = 1000.3333
$ws.Range("K45").Value = 857.3333
$ws.Range("L45").Value = 1000.3333
$ws.Range("M45").Value = -480.3333
$ws.Range("N45").Value = -1754.3333
$ws.Range("H61").Value = 1489.119
$ws.Range("I61").Value = 992.44446
$ws.Range("J61").Value = 4469.1665
$ws.Range("K61").Value = 992.44446
$ws.Range("L61").Value = 4469.1665
$ws.Range("M61").Value = -780.44446
$ws.Range("N61").Value = -4893.1665
$ws.Range("H74").Value = 4190.558
$ws.Range("I74").Value = 1257.3334
$ws.Range("J74").Value = 10959.538
$ws.Range("K74").Value = 1257.3334
$ws.Range("L74").Value = 10959.538
$ws.Range("M74").Value = -383.3334
$ws.Range("N74").Value = -12707.538
$ws.Range("H77").Value = 4190.558
$ws.Range("I77").Value = 1257.3334
$ws.Range("J77").Value = 10959.538
$ws.Range("K77").Value = 6286.666999999999
$ws.Range("L77").Value = 54797.69
$ws.Range("M77").Value = -1918.666999999999
$ws.Range("N77").Value = -63533.69
$ws.Range("H132").Value = 2142.275
$ws.Range("I132").Value = 1545.0667
$ws.Range("J132").Value = 3933.9
$ws.Range("K132").Value = 4635.2001
$ws.Range("L132").Value = 11801.7
$ws.Range("M132").Value = -2105.2001
$ws.Range("N132").Value = -16861.7
$ws.Range("H133").Value = 49200
$ws.Range("J133").Value = 49200
$ws.Range("L133").Value = 49200
$ws.Range("N133").Value = -54260
$ws.Range("H136").Value = 1489.119
$ws.Range("I136").Value = 992.44446
$ws.Range("J136").Value = 4469.1665
$ws.Range("K136").Value = 2977.33338
$ws.Range("L136").Value = 13407.4995
$ws.Range("M136").Value = -427.33338
$ws.Range("N136").Value = -18507.4995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1456.2046
$ws.Range("I20").Value = 1515.1666
$ws.Range("J20").Value = 1329.8572
$ws.Range("K20").Value = 1515.1666
$ws.Range("L20").Value = 1329.8572
$ws.Range("M20").Value = -1268.1666
$ws.Range("N20").Value = -1823.8572
$ws.Range("H134").Value = 17546430
$ws.Range("I134").Value = 25642460
$ws.Range("J134").Value = 5034.1113
$ws.Range("K134").Value = 76927380
$ws.Range("L134").Value = 15102.3339
$ws.Range("M134").Value = -76924845
$ws.Range("N134").Value = -20172.3339

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1081.5
$ws.Range("I22").Value = 463.55554
$ws.Range("J22").Value = 1876
$ws.Range("K22").Value = 463.55554
$ws.Range("L22").Value = 1876
$ws.Range("M22").Value = -113.55554
$ws.Range("N22").Value = -2576
$ws.Range("H31").Value = 1451.0217
$ws.Range("I31").Value = 934.7941
$ws.Range("K31").Value = 934.7941
$ws.Range("M31").Value = -639.7941
$ws.Range("H34").Value = 1451.0217
$ws.Range("I34").Value = 934.7941
$ws.Range("K34").Value = 934.7941
$ws.Range("M34").Value = -732.7941
$ws.Range("H48").Value = 10062
$ws.Range("J48").Value = 10062
$ws.Range("L48").Value = 10062
$ws.Range("N48").Value = -11014
$ws.Range("H58").Value = 1218.921
$ws.Range("I58").Value = 496.92593
$ws.Range("J58").Value = 2991.0908
$ws.Range("K58").Value = 496.92593
$ws.Range("L58").Value = 2991.0908
$ws.Range("M58").Value = -293.92593
$ws.Range("N58").Value = -3397.0908
$ws.Range("H107").Value = 271.12
$ws.Range("I107").Value = 188.375
$ws.Range("J107").Value = 310.05884
$ws.Range("K107").Value = 188.375
$ws.Range("L107").Value = 310.05884
$ws.Range("M107").Value = 1731.625
$ws.Range("N107").Value = -4150.05884
$ws.Range("H122").Value = 1762.4348
$ws.Range("I122").Value = 1008.4375
$ws.Range("J122").Value = 3485.8572
$ws.Range("K122").Value = 3025.3125
$ws.Range("L122").Value = 10457.5716
$ws.Range("M122").Value = -575.3125
$ws.Range("N122").Value = -15357.5716
$ws.Range("I132").Value = 2520.4285
$ws.Range("J132").Value = 3955
$ws.Range("K132").Value = 7561.2855
$ws.Range("L132").Value = 11865
$ws.Range("M132").Value = -5031.2855
$ws.Range("N132").Value = -16925
$ws.Range("H134").Value = 2940.7354
$ws.Range("I134").Value = 1578.5
$ws.Range("J134").Value = 6210.1
$ws.Range("K134").Value = 4735.5
$ws.Range("L134").Value = 18630.3
$ws.Range("M134").Value = -2200.5
$ws.Range("N134").Value = -23700.3
$ws.Range("H136").Value = 1218.921
$ws.Range("I136").Value = 496.92593
$ws.Range("J136").Value = 2991.0908
$ws.Range("K136").Value = 1490.77779
$ws.Range("L136").Value = 8973.2724
$ws.Range("M136").Value = 1059.22221
$ws.Range("N136").Value = -14073.2724

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 71428800
$ws.Range("I2").Value = 35.3
$ws.Range("J2").Value = 250000720
$ws.Range("K2").Value = 211.8
$ws.Range("L2").Value = 1500004320
$ws.Range("M2").Value = -98.79999999999998
$ws.Range("N2").Value = -1500004546
$ws.Range("H5").Value = 1192.7142
$ws.Range("I5").Value = 468.44446
$ws.Range("K5").Value = 1405.33338
$ws.Range("M5").Value = -1293.33338
$ws.Range("H135").Value = 1192.7142
$ws.Range("I135").Value = 468.44446
$ws.Range("K135").Value = 4216.00014
$ws.Range("M135").Value = -1681.00014

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5880.9116
$ws.Range("I70").Value = 6049.68
$ws.Range("J70").Value = 5412.1113
$ws.Range("K70").Value = 6049.68
$ws.Range("L70").Value = 5412.1113
$ws.Range("M70").Value = -5779.68
$ws.Range("N70").Value = -5952.1113
$ws.Range("H73").Value = 5880.9116
$ws.Range("I73").Value = 6049.68
$ws.Range("J73").Value = 5412.1113
$ws.Range("K73").Value = 6049.68
$ws.Range("L73").Value = 5412.1113
$ws.Range("M73").Value = -5113.68
$ws.Range("N73").Value = -7284.1113
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7084.8237
$ws.Range("I22").Value = 460.4
$ws.Range("J22").Value = 9845
$ws.Range("K22").Value = 460.4
$ws.Range("L22").Value = 9845
$ws.Range("M22").Value = -165.4
$ws.Range("N22").Value = -10435
$ws.Range("H27").Value = 7084.8237
$ws.Range("I27").Value = 460.4
$ws.Range("J27").Value = 9845
$ws.Range("K27").Value = 460.4
$ws.Range("L27").Value = 9845
$ws.Range("M27").Value = -353.4
$ws.Range("N27").Value = -10059
$ws.Range("H43").Value = 50000
$ws.Range("I43").Value = 50000
$ws.Range("K43").Value = 50000
$ws.Range("M43").Value = -49807

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H107").Value = 4274908.5
$ws.Range("I107").Value = 6945720
$ws.Range("J107").Value = 1610
$ws.Range("K107").Value = 20837160
$ws.Range("L107").Value = 4830
$ws.Range("M107").Value = -20835240
$ws.Range("N107").Value = -8670
$ws.Range("H122").Value = 39633.08
$ws.Range("I122").Value = 42854.793
$ws.Range("J122").Value = 972.5
$ws.Range("K122").Value = 128564.379
$ws.Range("L122").Value = 2917.5
$ws.Range("M122").Value = -126114.379
$ws.Range("N122").Value = -7817.5
$ws.Range("H132").Value = 17859688
$ws.Range("I132").Value = 21741150
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 65223450
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -65220920
$ws.Range("N132").Value = -19964
$ws.Range("H136").Value = 7269254.5
$ws.Range("I136").Value = 8155277
$ws.Range("J136").Value = 3869.4
$ws.Range("K136").Value = 24465831
$ws.Range("L136").Value = 11608.2
$ws.Range("M136").Value = -24463281
$ws.Range("N136").Value = -16708.2
